$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 617 ("「アラビア語が好き」" entry). This shifts every
# subsequent row up by one (618->617, 619->618, ..., 787->786) and the
# sheet's used range shrinks from A1:C787 to A1:C786.
$ws.Rows.Item(617).Delete()
